{"js": "// Resume update:\n//  1. \"Tech Stack\" line: add \"C# |\" after \"Python | Java | C |\"\n//  2. \"Frameworks and Libraries\" line: add \".NET |\" after \"Handlebars | Jinja2 |\"\n//  3. Remove the whole \"Leadership: ...\" bullet paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet techStackPara = null;\nlet frameworksPara = null;\nlet leadershipPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Tech Stack:\") === 0) {\n    techStackPara = p;\n  } else if (t.indexOf(\"Frameworks and Libraries:\") === 0) {\n    frameworksPara = p;\n  } else if (t.indexOf(\"Leadership:\") === 0) {\n    leadershipPara = p;\n  }\n}\n\nif (!techStackPara || !frameworksPara || !leadershipPara) {\n  throw new Error(\"Could not locate one of the target paragraphs.\");\n}\n\n// 1. Insert \"C# |\" right after \"Python | Java | C |\" in the Tech Stack line.\nconst cSharpAnchor = techStackPara.search(\"Python | Java | C | \", { matchCase: true });\ncSharpAnchor.load(\"text\");\nawait context.sync();\ncSharpAnchor.items[0].insertText(\"C# | \", Word.InsertLocation.after);\nawait context.sync();\n\n// 2. Insert \".NET |\" right after \"Handlebars | Jinja2 |\" in the Frameworks line.\nconst dotNetAnchor = frameworksPara.search(\"Jinja2 | \", { matchCase: true });\ndotNetAnchor.load(\"text\");\nawait context.sync();\ndotNetAnchor.items[0].insertText(\".NET | \", Word.InsertLocation.after);\nawait context.sync();\n\n// 3. Delete the entire \"Leadership: ...\" paragraph (bullet item).\nleadershipPara.delete();\nawait context.sync();\n", "ps1": "# Resume update:\n#  1. \"Tech Stack\" line: add \"C# |\" after \"Python | Java | C |\"\n#  2. \"Frameworks and Libraries\" line: add \".NET |\" after \"Handlebars | Jinja2 |\"\n#  3. Remove the whole \"Leadership: ...\" bullet paragraph.\n\n$d = $word.ActiveDocument\n\n# 1. Insert \"C# |\" right after \"Python | Java | C |\" in the Tech Stack line.\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\"Python | Java | C | \")\nif ($found1) {\n    $r1.InsertAfter(\"C# | \")\n}\n\n# 2. Insert \".NET |\" right after \"Handlebars | Jinja2 |\" in the Frameworks line.\n$r2 = $d.Content\n$found2 = $r2.Find.Execute(\"Jinja2 | \")\nif ($found2) {\n    $r2.InsertAfter(\".NET | \")\n}\n\n# 3. Delete the entire \"Leadership: ...\" paragraph (bullet item).\n$count = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"Leadership:\")) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -gt 0) {\n    $target = $d.Paragraphs.Item($targetIndex)\n    $target.Range.Delete()\n}\n"}
